$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates: force Text number format so the
# numeric-looking strings keep their exact original formatting
# (trailing zeros, leading zeros, precision) instead of Excel
# auto-converting them to native numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '242.56'
$ws.Range('D3').Value = '23.02'
$ws.Range('D4').Value = '5.422'
$ws.Range('D5').Value = '0.05887'
$ws.Range('D6').Value = '3.441'
$ws.Range('D7').Value = '6.537'
$ws.Range('D8').Value = '0.8105'
$ws.Range('D9').Value = '0.9412'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '0.1423'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '0.07436'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').Value = '0.03286'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '0.03052'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '0.09335'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').Value = '3.860'
$ws.Range('E15').Value = '14MCDexMCB'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = '0.001572'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = '0.04668'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').Value = '0.0005944'
$ws.Range('E18').Value = '17OneONE'
$ws.Range('D19').Value = '0.005886'
$ws.Range('D20').Value = '0.001261'
$ws.Range('D21').Value = '0.004884'
$ws.Range('D22').Value = '0.00006804'
$ws.Range('D23').Value = '3.573'
$ws.Range('D24').Value = '2.116'
$ws.Range('D26').Value = '0.1330'
$ws.Range('D27').Value = '0.0002286'
$ws.Range('D40').Value = '0.03948'
$ws.Range('D41').Value = '0.006190'
$ws.Range('D42').Value = '0.1069'
$ws.Range('D43').Value = '0.002571'
$ws.Range('D44').Value = '0.009334'
$ws.Range('E44').Value = '43LocalTradersLCTBestin24h'
$ws.Range('D45').Value = '0.00005181'
$ws.Range('D46').Value = '0.00000000751'
$ws.Range('D47').Value = '0.6785'
